$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: new B10 value; E10 adopts B-column's number-format/font style (stays empty) ---
$ws.Range("B10").Value = 0.00000000000001428
$ws.Range("B10").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 16: new B16 value; E16 adopts B-column's number-format/font style (stays empty) ---
$ws.Range("B16").Value = 50821.0752
$ws.Range("B16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 70: shift the old B70 value+style into E70, give B70 a new value, ---
# --- and add a brand-new (default-styled) value in D70                    ---
$ws.Range("B70").Copy($ws.Range("E70"))
$excel.CutCopyMode = $false
$ws.Range("B70").Value = 0.00000000000000006615
$ws.Range("D70").Value = 1.5

# --- Update which cell is shown as selected when the workbook is reopened ---
[void]$ws.Range("G17").Select()

# --- Rename the tab, and roll the sheet's internal id forward by duplicating ---
# --- it and dropping the original (mirrors re-saving the sheet under a new name) ---
[void]$ws.Copy($null, $ws)
[void]$wb.Worksheets.Item(1).Delete()
$ws2 = $wb.Worksheets.Item(1)
$ws2.Name = "230930"
[void]$ws2.Select()

Write-Host "Applied K-values update."
